$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 9).Value = 'Other found locations'

$ws.Cells.Item(2, 5).Value = '[Amaro N%Duarte‐Neto%NULL%1,  Renata A A%Monteiro%NULL%2,  Renata A A%Monteiro%NULL%0,  Luiz F F%da Silva%NULL%1,  Denise M A C%Malheiros%NULL%1,  Ellen P%de Oliveira%NULL%1,  Jair%Theodoro‐Filho%NULL%1,  João R R%Pinho%NULL%1,  Michele S%Gomes‐Gouvêa%NULL%1,  Ana P M%Salles%NULL%1,  Ilka R S%de Oliveira%NULL%1,  Thais%Mauad%NULL%1,  Paulo H N%Saldiva%NULL%1,  Marisa%Dolhnikoff%maridol@usp.br%1]'
$ws.Cells.Item(2, 9).Value = '_PMC'

$ws.Cells.Item(3, 9).Value = ''

$ws.Cells.Item(4, 5).Value = '[Louis Maximilian%Buja%l.maximilian.buja@uth.tmc.edu%1,  Dwayne A.%Wolf%NULL%1,  Bihong%Zhao%NULL%1,  Bindu%Akkanti%NULL%1,  Michelle%McDonald%NULL%1,  Laura%Lelenwa%NULL%1,  Noah%Reilly%NULL%1,  Giulia%Ottaviani%NULL%1,  M. Tarek%Elghetany%NULL%1,  Daniel Ocazionez%Trujillo%NULL%1,  Gabriel M.%Aisenberg%NULL%1,  Mohammad%Madjid%NULL%1,  Biswajit%Kar%NULL%1]'
$ws.Cells.Item(4, 9).Value = '_PMC_elsevier'

$ws.Cells.Item(5, 9).Value = ''

$ws.Cells.Item(6, 5).Value = '[J. Matthew%Lacy%NULL%1,  Erin G.%Brooks%NULL%1,  Joshua%Akers%NULL%1,  Danielle%Armstrong%NULL%1,  Lauren%Decker%NULL%1,  Adam%Gonzalez%NULL%1,  William%Humphrey%NULL%1,  Romana%Mayer%NULL%1,  Matthew%Miller%NULL%1,  Catherine%Perez%NULL%1,  Jose Antonio Ruiz%Arango%NULL%1,  Lakshmanan%Sathyavagiswaran%NULL%1,  Wendy%Stroh%NULL%1,  Suzanne%Utley%NULL%1]'
$ws.Cells.Item(6, 9).Value = '_PMC'

$ws.Cells.Item(7, 5).Value = '[Dominic%Wichmann%NULL%1,  Jan-Peter%Sperhake%NULL%2,  Marc%Lütgehetmann%NULL%2,  Stefan%Steurer%NULL%2,  Carolin%Edler%NULL%1,  Axel%Heinemann%NULL%2,  Fabian%Heinrich%NULL%2,  Herbert%Mushumba%NULL%2,  Inga%Kniep%NULL%1,  Ann Sophie%Schröder%NULL%2,  Christoph%Burdelski%NULL%2,  Geraldine%de Heer%NULL%2,  Axel%Nierhaus%NULL%2,  Daniel%Frings%NULL%2,  Susanne%Pfefferle%NULL%1,  Heinrich%Becker%NULL%1,  Hanns%Bredereke-Wiedling%NULL%1,  Andreas%de Weerth%NULL%2,  Hans-Richard%Paschen%NULL%1,  Sara%Sheikhzadeh-Eggers%NULL%1,  Axel%Stang%NULL%1,  Stefan%Schmiedel%NULL%1,  Carsten%Bokemeyer%NULL%1,  Marylyn M.%Addo%NULL%1,  Martin%Aepfelbacher%NULL%2,  Klaus%Püschel%NULL%2,  Stefan%Kluge%NULL%3]'
$ws.Cells.Item(7, 9).Value = '_PMC'

$ws.Cells.Item(8, 5).Value = '[Thomas%Menter%NULL%2,  Jasmin D%Haslbauer%NULL%4,  Jasmin D%Haslbauer%NULL%0,  Ronny%Nienhold%NULL%2,  Spasenija%Savic%NULL%2,  Helmut%Hopfer%NULL%2,  Nikolaus%Deigendesch%NULL%2,  Stephan%Frank%NULL%2,  Daniel%Turek%NULL%2,  Niels%Willi%NULL%2,  Hans%Pargger%NULL%2,  Stefano%Bassetti%NULL%2,  Joerg D%Leuppi%NULL%2,  Gieri%Cathomas%NULL%2,  Markus%Tolnay%NULL%2,  Kirsten D%Mertz%NULL%2,  Alexandar%Tzankov%alexandar.tzankov@usb.ch%2]'
$ws.Cells.Item(8, 9).Value = '_PMC'

$ws.Cells.Item(9, 5).Value = '[Thomas%Menter%NULL%0,  Jasmin D%Haslbauer%NULL%0,  Jasmin D%Haslbauer%NULL%0,  Ronny%Nienhold%NULL%0,  Spasenija%Savic%NULL%0,  Helmut%Hopfer%NULL%0,  Nikolaus%Deigendesch%NULL%0,  Stephan%Frank%NULL%0,  Daniel%Turek%NULL%0,  Niels%Willi%NULL%0,  Hans%Pargger%NULL%0,  Stefano%Bassetti%NULL%0,  Joerg D%Leuppi%NULL%0,  Gieri%Cathomas%NULL%0,  Markus%Tolnay%NULL%0,  Kirsten D%Mertz%NULL%0,  Alexandar%Tzankov%alexandar.tzankov@usb.ch%0]'
$ws.Cells.Item(9, 9).Value = '_PMC'

$ws.Cells.Item(10, 5).Value = '[Zsuzsanna%Varga%NULL%1,  Andreas J%Flammer%NULL%1,  Peter%Steiger%NULL%1,  Martina%Haberecker%NULL%1,  Rea%Andermatt%NULL%1,  Annelies S%Zinkernagel%NULL%1,  Mandeep R%Mehra%NULL%1,  Reto A%Schuepbach%NULL%1,  Frank%Ruschitzka%frank.ruschitzka@usz.ch%1,  Holger%Moch%NULL%1]'
$ws.Cells.Item(10, 9).Value = '_PMC_elsevier'

$ws.Cells.Item(11, 5).Value = '[Sufang%Tian%NULL%1,  Yong%Xiong%NULL%0,  Huan%Liu%NULL%3,  Li%Niu%NULL%1,  Jianchun%Guo%NULL%1,  Meiyan%Liao%NULL%1,  Shu-Yuan%Xiao%NULL%1]'
$ws.Cells.Item(11, 9).Value = '_PMC_elsevier_Springer'

$ws.Cells.Item(12, 5).Value = '[Lisa M%Barton%NULL%1,  Eric J%Duval%NULL%1,  Edana%Stroberg%NULL%1,  Subha%Ghosh%NULL%1,  Sanjay%Mukhopadhyay%mukhops@ccf.org%1]'
$ws.Cells.Item(12, 9).Value = '_PMC'

$ws.Cells.Item(13, 5).Value = '[Pedro%Navarro Conde%navarro_ped@gva.es%1,  Pilar%Alemany Monraval%NULL%1,  Carmen%Medina Medina%NULL%1,  Ana%Jiménez Sánchez%NULL%1,  Juan Carlos%Andrés Teruel%NULL%1,  José%Ferrando Marco%NULL%1,  Víctor%Puglia Santos%NULL%1,  Empar%Mayordomo Aranda%NULL%1]'
$ws.Cells.Item(13, 9).Value = '_PMC_elsevier'

$ws.Cells.Item(14, 5).Value = '[Carolin%Edler%NULL%1,  Ann Sophie%Schröder%NULL%0,  Martin%Aepfelbacher%NULL%0,  Antonia%Fitzek%NULL%1,  Axel%Heinemann%NULL%0,  Fabian%Heinrich%NULL%0,  Anke%Klein%NULL%1,  Felicia%Langenwalder%NULL%1,  Marc%Lütgehetmann%NULL%0,  Kira%Meißner%NULL%1,  Klaus%Püschel%NULL%0,  Julia%Schädler%NULL%1,  Stefan%Steurer%NULL%0,  Herbert%Mushumba%NULL%0,  Jan-Peter%Sperhake%sperhake@uke.de%0]'
$ws.Cells.Item(14, 9).Value = '_PMC_Springer'

$ws.Cells.Item(15, 9).Value = ''

$ws.Cells.Item(16, 5).Value = '[Miroslav%Sekulic%NULL%1,  Holly%Harper%NULL%1,  Behtash G%Nezami%NULL%1,  Daniel L%Shen%NULL%1,  Simona Pichler%Sekulic%NULL%1,  Aaron T%Koeth%NULL%1,  Clifford V%Harding%NULL%1,  Hannah%Gilmore%NULL%1,  Navid%Sadri%navid.sadri@uhhospitals.org%1]'
$ws.Cells.Item(16, 9).Value = '_PMC'

$ws.Cells.Item(17, 5).Value = '[Christine%Suess%NULL%1,  Roland%Hausmann%roland.hausmann@kssg.ch%1]'
$ws.Cells.Item(17, 9).Value = '_PMC_Springer'

$ws.Cells.Item(18, 5).Value = '[Diego%Aguiar%diego.aguiar@hcuge.ch%1,  Johannes Alexander%Lobrinus%NULL%2,  Johannes Alexander%Lobrinus%NULL%0,  Manuel%Schibler%NULL%1,  Tony%Fracasso%NULL%1,  Christelle%Lardi%NULL%1]'
$ws.Cells.Item(18, 9).Value = '_PMC_Springer'

$ws.Cells.Item(19, 5).Value = '[Sharon E%Fox%NULL%1,  Aibek%Akmatbekov%NULL%1,  Jack L%Harbert%NULL%1,  Guang%Li%NULL%1,  J%Quincy Brown%NULL%1,  Richard S%Vander Heide%NULL%1]'
$ws.Cells.Item(19, 9).Value = '_PMC_elsevier'

$ws.Cells.Item(20, 5).Value = '[Mohammad Taghi%Beigmohammadi%NULL%1,  Behnaz%Jahanbin%NULL%1,  Masoomeh%Safaei%NULL%1,  Laya%Amoozadeh%NULL%1,  Meysam%Khoshavi%NULL%1,  Vahid%Mehrtash%NULL%1,  Bita%Jafarzadeh%NULL%1,  Alireza%Abdollahi%NULL%1]'
$ws.Cells.Item(20, 9).Value = '_PMC'

$ws.Cells.Item(21, 5).Value = '[Chaofu%Wang%wangchaofu@126.com%1,  Jing%Xie%NULL%2,  Lei%Zhao%NULL%1,  Xiaochun%Fei%NULL%1,  Heng%Zhang%NULL%1,  Yun%Tan%NULL%1,  Xiu%Nie%NULL%1,  Luting%Zhou%NULL%1,  Zhenhua%Liu%NULL%1,  Yong%Ren%NULL%1,  Ling%Yuan%NULL%1,  Yu%Zhang%NULL%1,  Jinsheng%Zhang%NULL%1,  Liwei%Liang%NULL%1,  Xinwei%Chen%NULL%1,  Xin%Liu%NULL%2,  Peng%Wang%NULL%1,  Xiao%Han%NULL%1,  Xiangqin%Weng%NULL%1,  Ying%Chen%NULL%1,  Ting%Yu%NULL%0,  Xinxin%Zhang%NULL%0,  Jun%Cai%caijun@shsmu.edu.cn%1,  Rong%Chen%crjudy@126.com%3,  Zhengli%Shi%zlshi@wh.iov.cn%2,  Xiuwu%Bian%bianxiuwu@263.net%1]'
$ws.Cells.Item(21, 9).Value = '_PMC_elsevier'

$ws.Cells.Item(22, 5).Value = '[Hans%Bösmüller%hans.boesmueller@med.uni-tuebingen.de%1,  Selina%Traxler%NULL%2,  Selina%Traxler%NULL%0,  Michael%Bitzer%NULL%1,  Helene%Häberle%NULL%1,  Wolfgang%Raiser%NULL%1,  Dominik%Nann%NULL%1,  Leonie%Frauenfeld%NULL%1,  Antonio%Vogelsberg%NULL%1,  Karin%Klingel%NULL%1,  Falko%Fend%NULL%1]'
$ws.Cells.Item(22, 9).Value = '_PMC_Springer'

$ws.Cells.Item(23, 5).Value = '[Wolf%Schweitzer%wolf.schweitzer@irm.uzh.ch%1,  Thomas%Ruder%NULL%1,  Rilana%Baumeister%NULL%1,  Stephan%Bolliger%NULL%1,  Michael%Thali%NULL%1,  Eva%Meixner%NULL%1,  Garyfalia%Ampanozi%NULL%1]'
$ws.Cells.Item(23, 9).Value = '_PMC_elsevier'

$ws.Cells.Item(24, 5).Value = '[Zhe%Xu%NULL%1,  Lei%Shi%NULL%1,  Yijin%Wang%NULL%1,  Jiyuan%Zhang%NULL%1,  Lei%Huang%NULL%2,  Chao%Zhang%NULL%1,  Shuhong%Liu%NULL%1,  Peng%Zhao%NULL%1,  Hongxia%Liu%NULL%1,  Li%Zhu%NULL%1,  Yanhong%Tai%NULL%1,  Changqing%Bai%NULL%1,  Tingting%Gao%NULL%1,  Jinwen%Song%NULL%1,  Peng%Xia%NULL%1,  Jinghui%Dong%NULL%1,  Jingmin%Zhao%NULL%1,  Fu-Sheng%Wang%NULL%1]'
$ws.Cells.Item(24, 9).Value = '_PMC_elsevier'

$ws.Cells.Item(25, 9).Value = ''

$ws.Cells.Item(26, 5).Value = '[Mathilde%Ducloyer%Mathilde.ducloyer@chu-nantes.fr%1,  Benjamin%Gaborit%NULL%2,  Benjamin%Gaborit%NULL%0,  Claire%Toquet%NULL%1,  Louise%Castain%NULL%1,  Antonin%Bal%NULL%1,  Pierre Paul%Arrigoni%NULL%1,  Raphaël%Lecomte%NULL%1,  Renaud%Clement%NULL%1,  Christine%Sagan%NULL%1]'
$ws.Cells.Item(26, 9).Value = '_PMC_Springer'

$ws.Cells.Item(27, 5).Value = '[Andreea-Elena%Cîrstea%NULL%1,  Radu Lucian%Buzulică%NULL%1,  Daniel%Pirici%NULL%1,  Mihail Constantin%Ceauşu%NULL%1,  Radu Vasile%Iman%NULL%1,  Ovidiu-Mircea%Gheorghe%NULL%1,  Simona Daniela%Neamţu%NULL%1,  Liliana%Stanca%NULL%1,  Răzvan%Ene%NULL%1,  Samir%Kumar-Singh%NULL%1,  Laurenţiu%Mogoantă%NULL%1]'
$ws.Cells.Item(27, 9).Value = '_PMC'

$ws.Cells.Item(28, 5).Value = '[Leo%Nicolai%NULL%1,  Alexander%Leunig%NULL%1,  Sophia%Brambs%NULL%1,  Rainer%Kaiser%NULL%1,  Tobias%Weinberger%NULL%1,  Michael%Weigand%NULL%1,  Maximilian%Muenchhoff%NULL%1,  Johannes C.%Hellmuth%NULL%1,  Stephan%Ledderose%NULL%1,  Heiko%Schulz%NULL%1,  Clemens%Scherer%NULL%1,  Martina%Rudelius%NULL%1,  Michael%Zoller%NULL%1,  Dominik%Höchter%NULL%1,  Oliver%Keppler%NULL%1,  Daniel%Teupser%NULL%1,  Bernhard%Zwißler%NULL%1,  Michael%von Bergwelt-Baildon%NULL%1,  Stefan%Kääb%NULL%1,  Steffen%Massberg%NULL%1,  Kami%Pekayvaz%NULL%1,  Konstantin%Stark%NULL%1]'
$ws.Cells.Item(28, 9).Value = '_PMC'

$ws.Cells.Item(29, 5).Value = '[Claudia%Grosse%NULL%1,  Alexandra%Grosse%NULL%1,  Helmut J.F.%Salzer%NULL%1,  Martin W.%Dünser%NULL%1,  Reinhard%Motz%NULL%1,  Rupert%Langer%NULL%1]'
$ws.Cells.Item(29, 9).Value = '_PMC_elsevier'

$ws.Cells.Item(30, 9).Value = ''

$ws.Cells.Item(31, 5).Value = '[Myriam%Remmelink%NULL%1,  Ricardo%De Mendonça%NULL%1,  Nicky%D’Haene%NULL%1,  Sarah%De Clercq%NULL%1,  Camille%Verocq%NULL%1,  Laetitia%Lebrun%NULL%1,  Philomène%Lavis%NULL%1,  Marie-Lucie%Racu%NULL%1,  Anne-Laure%Trépant%NULL%1,  Calliope%Maris%NULL%1,  Sandrine%Rorive%NULL%1,  Jean-Christophe%Goffard%NULL%1,  Olivier%De Witte%NULL%1,  Lorenzo%Peluso%NULL%1,  Jean-Louis%Vincent%NULL%1,  Christine%Decaestecker%NULL%1,  Fabio Silvio%Taccone%NULL%1,  Isabelle%Salmon%Isabelle.Salmon@erasme.ulb.ac.be%1]'
$ws.Cells.Item(31, 9).Value = '_PMC_Springer'

$ws.Cells.Item(32, 5).Value = '[Koji%Okudela%kojixok@yokohama.cu.ac.jp%1,  Hiroyuki%Hayashi%NULL%2,  Hiroyuki%Hayashi%NULL%0,  Yukihiro%Yoshimura%NULL%1,  Hiroaki%Sasaki%NULL%1,  Hiroshi%Horiuchi%NULL%1,  Nobuyuki%Miyata%NULL%1,  Natsuo%Tachikawa%NULL%1,  Yuki%Tsuchiya%NULL%1,  Hideaki%Mitsui%NULL%1,  Kenichi%Ohashi%NULL%1]'
$ws.Cells.Item(32, 9).Value = '_PMC'

$ws.Cells.Item(33, 5).Value = '[Takuya%Adachi%NULL%1,  Ja-Mun%Chong%NULL%1,  Noriko%Nakajima%NULL%1,  Masahiro%Sano%NULL%1,  Jun%Yamazaki%NULL%1,  Ippei%Miyamoto%NULL%1,  Haruka%Nishioka%NULL%1,  Hidetaka%Akita%NULL%1,  Yuko%Sato%NULL%1,  Michiyo%Kataoka%NULL%1,  Harutaka%Katano%NULL%1,  Minoru%Tobiume%NULL%1,  Tsuyoshi%Sekizuka%NULL%1,  Kentaro%Itokawa%NULL%1,  Makoto%Kuroda%NULL%1,  Tadaki%Suzuki%NULL%1]'
$ws.Cells.Item(33, 9).Value = '_PMC'

$ws.Cells.Item(34, 5).Value = '[Girish N.%Nadkarni%NULL%1,  Anuradha%Lala%NULL%1,  Emilia%Bagiella%NULL%1,  Helena L.%Chang%NULL%1,  Pedro R.%Moreno%NULL%1,  Elisabet%Pujadas%NULL%1,  Varun%Arvind%NULL%1,  Sonali%Bose%NULL%1,  Alexander W.%Charney%NULL%1,  Martin D.%Chen%NULL%1,  Carlos%Cordon-Cardo%NULL%1,  Andrew S.%Dunn%NULL%1,  Michael E.%Farkouh%NULL%1,  Benjamin S.%Glicksberg%NULL%1,  Arash%Kia%NULL%1,  Roopa%Kohli-Seth%NULL%1,  Matthew A.%Levin%NULL%1,  Prem%Timsina%NULL%1,  Shan%Zhao%NULL%1,  Zahi A.%Fayad%NULL%5,  Valentin%Fuster%NULL%1]'
$ws.Cells.Item(34, 9).Value = '_PMC_elsevier'

$ws.Cells.Item(35, 5).Value = '[Osama%Al-Dalahmah%oa2298@cumc.columbia.edu%1,  Kiran T.%Thakur%NULL%2,  Kiran T.%Thakur%NULL%0,  Anna S.%Nordvig%NULL%1,  Morgan L.%Prust%NULL%1,  William%Roth%NULL%1,  Angela%Lignelli%NULL%1,  Anne-Catrin%Uhlemann%NULL%1,  Emily Happy%Miller%NULL%1,  Shajo%Kunnath-Velayudhan%NULL%1,  Armando%Del Portillo%NULL%1,  Yang%Liu%NULL%1,  Gunnar%Hargus%NULL%1,  Andrew F.%Teich%NULL%1,  Richard A.%Hickman%NULL%1,  Kurenai%Tanji%NULL%1,  James E.%Goldman%NULL%1,  Phyllis L.%Faust%NULL%1,  Peter%Canoll%NULL%1]'
$ws.Cells.Item(35, 9).Value = '_PMC_Springer'

$ws.Cells.Item(36, 5).Value = '[George-Călin%Oprinca%george_oprinca@yahoo.com%1,  Lilioara-Alexandra%Muja%NULL%2,  Lilioara-Alexandra%Muja%NULL%0]'
$ws.Cells.Item(36, 9).Value = '_PMC_Springer'

$ws.Cells.Item(37, 5).Value = '[Xin-xin%Wang%NULL%1,  Chen%Shao%NULL%2,  Chen%Shao%NULL%0,  Xiao-jie%Huang%NULL%1,  Lin%Sun%NULL%1,  Ling-jia%Meng%NULL%1,  Hui%Liu%NULL%2,  Shi-jie%Zhang%NULL%1,  Hong-jun%Li%NULL%1,  Fu-dong%Lv%NULL%1]'
$ws.Cells.Item(37, 9).Value = '_PMC'

$ws.Cells.Item(38, 9).Value = ''

$ws.Cells.Item(39, 5).Value = '[Sarah S.%Elsoukkary%NULL%1,  Maria%Mostyka%NULL%1,  Alicia%Dillard%NULL%1,  Diana R.%Berman%NULL%1,  Lucy X.%Ma%NULL%1,  Amy%Chadburn%NULL%1,  Rhonda K.%Yantiss%NULL%1,  Jose%Jessurun%NULL%1,  Surya V.%Seshan%NULL%1,  Alain C.%Borczuk%NULL%1,  Steven P.%Salvatore%NULL%1]'
$ws.Cells.Item(39, 9).Value = '_PMC'

$ws.Cells.Item(40, 5).Value = '[Brian%Hanley%NULL%1,  Kikkeri N%Naresh%NULL%1,  Candice%Roufosse%NULL%1,  Andrew G%Nicholson%NULL%1,  Justin%Weir%NULL%1,  Graham S%Cooke%NULL%1,  Mark%Thursz%NULL%1,  Pinelopi%Manousou%NULL%1,  Richard%Corbett%NULL%1,  Robert%Goldin%NULL%1,  Safa%Al-Sarraj%NULL%1,  Alireza%Abdolrasouli%NULL%1,  Olivia C%Swann%NULL%1,  Laury%Baillon%NULL%1,  Rebecca%Penn%NULL%1,  Wendy S%Barclay%NULL%1,  Patrizia%Viola%NULL%1,  Michael%Osborn%NULL%1]'
$ws.Cells.Item(40, 9).Value = '_PMC_elsevier'
